$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update existing rows 12-18 (renaming "Main 3,3Buck ..." designators to the new Buck schema) ---

# Row 12: TPS54308 U1 -> qty 1 to 2, designator renamed
$ws.Cells.Item(12, 1).Value = 2
$ws.Cells.Item(12, 6).Value = "Main Buck U1"

# Row 13: 10uF cap Cin -> qty 1 to 2, designator renamed
$ws.Cells.Item(13, 1).Value = 2
$ws.Cells.Item(13, 6).Value = "Main Buck Cin"

# Row 14: Panasonic 33uF cap Cout -> qty 1 to 2, designator renamed
$ws.Cells.Item(14, 1).Value = 2
$ws.Cells.Item(14, 6).Value = "Main Buck Cout"

# Row 15: Wurth inductor -> designator renamed (L1 -> 3.3Buck L2)
$ws.Cells.Item(15, 6).Value = "Main 3.3Buck L2"

# Row 16: SMD Resistor 0805 22kohm -> designator renamed (Rfbb -> 3.3Buck Rfbb2)
$ws.Cells.Item(16, 6).Value = "Main 3.3Buck Rfbb2"

# Row 17: SMD Chip Resistor 0603 100kohm -> designator renamed (Rfbt -> Buck Rfbt)
$ws.Cells.Item(17, 6).Value = "Main Buck Rfbt"

# Row 18: SMD Multilayer Ceramic Capacitor 0.1uF -> designator renamed (Cbst -> Buck Cbst)
$ws.Cells.Item(18, 6).Value = "Main Buck Cbst"

# --- Add new rows 19-20 ---

$ws.Cells.Item(19, 1).Value = 1
$ws.Cells.Item(19, 2).Value = "7447798271 -  Power Inductor (SMD), 27 µH"
$ws.Cells.Item(19, 3).Value = 1800291
$ws.Cells.Item(19, 4).Value = "farnell"
$ws.Cells.Item(19, 5).Value = 3.77
$ws.Cells.Item(19, 6).Value = "Main 5Buck L1"

$ws.Cells.Item(20, 1).Value = 1
$ws.Cells.Item(20, 2).Value = "CPF0603F13K7C1 -  SMD Chip Resistor, 0603, 13.7 kohm"
$ws.Cells.Item(20, 3).Value = 1527596
$ws.Cells.Item(20, 4).Value = "farnell"
$ws.Cells.Item(20, 5).Value = 0.168
$ws.Cells.Item(20, 6).Value = "Main 5buck Rfbb1"

# --- Formatting / view updates ---

# Column F width adjustment (closest best-fit width to the new, longer designator text)
$ws.Columns.Item(6).ColumnWidth = 17.6

# Selection moved to F19
$ws.Range("F19").Select()
